# Changed Card Numbers for FeeData
#
# Sheet1 ("Account" column, col C) held several test card numbers that need
# to be swapped out:
#   - 4111111111111111  ->  4761730000000011   (rows 2-5)
#   - 4400000000000008  ->  4000056655665556   (row 8)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("C2").Value = "4761730000000011"
$ws1.Range("C3").Value = "4761730000000011"
$ws1.Range("C4").Value = "4761730000000011"
$ws1.Range("C5").Value = "4761730000000011"

$ws1.Range("C8").Value = "4000056655665556"

# Leave the workbook focused on Sheet1 / the last cell touched, matching
# where the editor ended up after making the change.
$ws1.Activate()
$ws1.Range("C8").Select()
